$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Reviewed by" (column M) notes added while reviewing some test cases:
# rows 2-5 get "Fatma" as the reviewer.
$ws.Range("M2").Value = "Fatma"
$ws.Range("M3").Value = "Fatma"
$ws.Range("M4").Value = "Fatma"
$ws.Range("M5").Value = "Fatma"

# Match the center/middle wrap alignment used elsewhere in the sheet for
# the "Bug ID" column cells (style shared with L2/L5/L9/L10/L12/L19/L25).
# (WrapText is already on for these cells, so only the alignment needs
# touching - this keeps the engine from minting extra unused cell formats.)
$bugIdCells = @("L2", "L5", "L9", "L10", "L12", "L19", "L25")
foreach ($addr in $bugIdCells) {
    $cell = $ws.Range($addr)
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# Move the cursor/selection to L3 (reviewer scrolled back up while checking
# the bug-id column after the edits).
$ws.Range("L3").Select() | Out-Null
